$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7515.90007656357
$ws.Range("C2").Value = 7013.36567061664
$ws.Range("E2").Value = 2845.52588361621
$ws.Range("F2").Value = 42.9538147597023

$ws.Range("B3").Value = 6871.50363258918
$ws.Range("C3").Value = 5788.42714292018
$ws.Range("E3").Value = 3196.11820704134
$ws.Range("F3").Value = 223.522722915063

$ws.Range("B4").Value = 2040.99989232529
$ws.Range("C4").Value = 3393.45998970999
$ws.Range("E4").Value = 3144.40012534721
$ws.Range("F4").Value = -102.755828539283

$ws.Range("B5").Value = 1979.57791610586
$ws.Range("C5").Value = 3553.35299879906
$ws.Range("E5").Value = 3071.37502233917
$ws.Range("F5").Value = -99.1363324525737

$ws.Range("B6").Value = 7935.63815207031
$ws.Range("C6").Value = 6836.83190900488
$ws.Range("E6").Value = 4179.18449116543
$ws.Range("F6").Value = 83.8340166737629

$ws.Range("B7").Value = 7961.86044293685
$ws.Range("C7").Value = 7458.0984026189
$ws.Range("E7").Value = 4204.66862491391
$ws.Range("F7").Value = 110.781959480534

$ws.Range("C9").Value = 8058.67158058859
$ws.Range("F9").Value = 151.809595379654

$ws.Range("C10").Value = 7207.15270225394
$ws.Range("F10").Value = 117.60581955344

$ws.Range("C11").Value = 5119.04966860091
$ws.Range("F11").Value = 17.6293133344348

$ws.Range("C12").Value = 5108.8029688938
$ws.Range("F12").Value = 16.8996130165795

$ws.Range("C13").Value = 8373.00230763601
$ws.Range("F13").Value = 192.356497479303

$ws.Range("C14").Value = 8643.03998400068
$ws.Range("F14").Value = 203.60822416554

$ws.Range("C15").Value = 8708.4495355909
$ws.Range("F15").Value = 206.333622148466
